# Multi browser testing code updated
# Work on the DATA worksheet (sheet2)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# 1. Insert a new column before column C. This shifts the existing
#    C (username) and D (password) columns to D and E respectively,
#    and shifts the row/col dimension + cell contents/styles along with it.
$ws.Columns("C").Insert()

# 2. Populate the new "browser" column and update the row-2 password.
#    The order of these assignments matters: it controls the order in
#    which brand-new shared strings are appended to the shared string
#    table, which must match: Asdf@1234, browser, chrome, firefox, microedge
$ws.Range("E2").Value = "Asdf@1234"
$ws.Range("C1").Value = "browser"
$ws.Range("C2").Value = "chrome"
$ws.Range("C6").Value = "firefox"
$ws.Range("C3").Value = "microedge"
$ws.Range("C4").Value = "chrome"
$ws.Range("C5").Value = "chrome"

# 3. The column insert does not re-point the existing mailto hyperlinks
#    at their new (shifted) cells, so rebuild the hyperlinks collection
#    from scratch pointing at the new E2:E6 cells, preserving the
#    original target addresses/order (rId1..rId5 keep the same targets).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Asdf@123")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:Asdf@123")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:Asdf@333")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:Asdf@123")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:Asdfg@123")

# Hyperlinks.Add() re-applies its own (slightly different) cell style;
# restore the original "Hyperlink" cell style so these cells keep using
# the same style index as before the edit.
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E6").Style = "Hyperlink"

# 4. The newly inserted column doesn't inherit column B's custom width;
#    set it to match (closest achievable width to 14.28515625).
$ws.Columns("C").ColumnWidth = 13.43

# 5. Update the remembered selection/active cell.
[void]$ws.Range("C15").Select()
